$d = $word.ActiveDocument

# --- Edit 1: "After running queries..." paragraph ---
$d.Content.Find.Execute('to "arl-tally.txt"   Be sure to ', $true, $false, $false, $false, $false, $true, 1, $false, 'to <arl-tally.txt> (Be sure to ', 2) | Out-Null
$d.Content.Find.Execute(' of the script. ', $true, $false, $false, $false, $false, $true, 1, $false, ' of the script).  Edit the outfile to format for Excel, and save as <arl_tally_formated4excel.txt>', 2) | Out-Null

# --- Edit 2: "Add the info..." paragraph ---
$d.Content.Find.Execute('Add the info from "arl-tally" to ', $true, $false, $false, $false, $false, $true, 1, $false, 'Copy and paste from <arl_tally_formated4excel.txt> to ', 2) | Out-Null
$d.Content.Find.Execute('spreadsheet', $true, $false, $false, $false, $false, $true, 1, $false, 'EXCEL Workbook, sheet 1 “yy-yy-q1 all formats”', 2) | Out-Null

# --- Edit 3: move _GoBack bookmark + insert new paragraph about query results ---
$p413 = $d.Paragraphs(413)
$rng = $p413.Range
$rng.InsertBefore("PLACEHOLDER_BM")
$p413b = $d.Paragraphs(413)
$bmRng = $p413b.Range
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng)
$d.Content.Find.Execute("PLACEHOLDER_BM", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$p413c = $d.Paragraphs(413)
$p413c.Range.InsertParagraphAfter()
$d.Paragraphs(414).Range.Text = "The query results include all University of Cincinnati Libraries (Univ. Libs, UCBA, Clermont, HSL & Law)"

# --- Edit 4: "711 bibs" paragraph ---
$d.Content.Find.Execute("there are only 711 bibs with ", $true, $false, $false, $false, $false, $true, 1, $false, "there are approximately 711 bibs with ", 2) | Out-Null

Write-Output "all edits done"
